$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three fixtures played on 2023-11-04 (row numbers 98-100, all sharing
# the same kickoff date) were re-ordered by the upstream data refresh.
# Rotate their data (every column except the row-sequence column A) so the
# row that used to be 100 becomes 98, the old 98 becomes 99, and the old 99
# becomes 100.
$oldRow98 = $ws.Range("B98:AB98").Value2
$oldRow99 = $ws.Range("B99:AB99").Value2
$oldRow100 = $ws.Range("B100:AB100").Value2

$ws.Range("B98:AB98").Value = $oldRow100
$ws.Range("B99:AB99").Value = $oldRow98
$ws.Range("B100:AB100").Value = $oldRow99

# Drop the placeholder fixture (id "7874808", Kaisar Kyzylorda vs Kairat
# Almaty) that had no final score yet - it was removed from the source feed.
$ws.Rows(136).Delete()
